$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'62.230.51"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -1.38%  '
$ws.Range("D3").Value = "'3.002.42"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.53%  '
$ws.Range("E4").Value = '  +0.04%  '
$ws.Range("D5").Value = "'586.58"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.11%  '
$ws.Range("D6").Value = "'146.21"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -3.32%  '
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("E8").Value = '  -2.20%  '
$ws.Range("D9").Value = "'2.999.43"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.66%  '
$ws.Range("D10").Value = "'0.147"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -4.11%  '
$ws.Range("D11").Value = "'5.76"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.57%  '
$ws.Range("E12").Value = '  +3.23%  '
$ws.Range("E13").Value = '  -2.41%  '
$ws.Range("D14").Value = "'34.49"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -4.91%  '
$ws.Range("E15").Value = '  +2.08%  '
$ws.Range("D16").Value = "'3.505.10"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -1.27%  '
$ws.Range("D17").Value = "'7.05"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -1.24%  '
$ws.Range("D18").Value = "'62.278.41"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -1.33%  '
$ws.Range("D19").Value = "'3.004.16"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -1.42%  '
$ws.Range("D20").Value = "'457.01"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  -4.02%  '
$ws.Range("D21").Value = "'13.96"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -2.12%  '
$ws.Range("D22").Value = "'0.687"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.46%  '
$ws.Range("D23").Value = "'7.39"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.57%  '
$ws.Range("D24").Value = "'81.69"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.60%  '
$ws.Range("E25").Value = '  -8.93%  '
$ws.Range("D26").Value = "'12.18"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -4.07%  '
$ws.Range("E27").Value = '  +0.01%  '
$ws.Range("D28").Value = "'9.78"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -7.93%  '
$ws.Range("E29").Value = '  -0.05%  '
$ws.Range("E30").Value = '  -1.28%  '
$ws.Range("D31").Value = "'6.94"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -5.48%  '
$ws.Range("D32").Value = "'2.08"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -5.09%  '
$ws.Range("D33").Value = "'27.68"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.20%  '
$ws.Range("E34").Value = '  -1.61%  '
$ws.Range("D35").Value = "'0.0₃0806"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -1.34%  '
$ws.Range("E36").Value = '  -2.92%  '
$ws.Range("D37").Value = "'5.74"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.34%  '
$ws.Range("E38").Value = '  -5.02%  '
$ws.Range("D39").Value = "'50.35"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -0.26%  '
$ws.Range("D40").Value = "'9.16"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.88%  '
$ws.Range("E41").Value = '  +7.79%  '
$ws.Range("D42").Value = "'2.88"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -10.97%  '
$ws.Range("D43").Value = "'393.36"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -9.14%  '
$ws.Range("E44").Value = '  -1.42%  '
$ws.Range("E45").Value = '  -7.10%  '
$ws.Range("D46").Value = "'2.728.61"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -3.44%  '
$ws.Range("D47").Value = "'37.55"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -1.91%  '
$ws.Range("D48").Value = "'129.18"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.39%  '
$ws.Range("E50").Value = '  -0.46%  '
$ws.Range("D51").Value = "'2.18"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.85%  '
